# Updates question numbers in every quiz sheet's result/question table as
# well as re-numbering the two "result" sheets (quiz2..quiz5) that used to
# mirror earlier quizzes (89-99, 64-88, 16-39, 40-63) so that every sheet's
# question numbers run sequentially 1-25 (quiz1 keeps 1-15 for rows 2-16 and
# continues 16-25 for rows 17-26). A new border/alignment style (thin box
# without a top edge, centered text) is applied to C3:C26 on quiz2..quiz5 to
# match the renumbered column. Finally the active sheet/tab moves from
# metadata_quiz to quiz6, and the stale view state (old selections / scroll
# position) on the quiz sheets is refreshed.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# quiz1: only the bottom half of the table (rows 17-26) was renumbered,
# continuing on from the sequential numbers already used in rows 2-16.
# No style/border change here, no selection change either.
# ---------------------------------------------------------------------
$quiz1 = $wb.Worksheets.Item("quiz1")
for ($r = 17; $r -le 26; $r++) {
    $quiz1.Cells.Item($r, 3).Value = ($r - 1)
}

# ---------------------------------------------------------------------
# quiz2, quiz3, quiz4, quiz5: renumber C2:C26 to 1..25, restyle C3:C26
# with a thin left/right/bottom border (no top) plus centered text, and
# move the stored selection to C2:C26 (clearing any stale scroll/selection
# state left over from editing).
# ---------------------------------------------------------------------
$resultSheets = @("quiz2", "quiz3", "quiz4", "quiz5")

foreach ($name in $resultSheets) {
    $ws = $wb.Worksheets.Item($name)

    for ($r = 2; $r -le 26; $r++) {
        $ws.Cells.Item($r, 3).Value = ($r - 1)
    }

    for ($r = 3; $r -le 26; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        $cell.Borders.Item(7).LineStyle = 1
        $cell.Borders.Item(7).Weight = 2
        $cell.Borders.Item(9).LineStyle = 1
        $cell.Borders.Item(9).Weight = 2
        $cell.Borders.Item(10).LineStyle = 1
        $cell.Borders.Item(10).Weight = 2
        $cell.Borders.Item(8).LineStyle = -4142
        $cell.HorizontalAlignment = -4108
        $cell.VerticalAlignment = -4108
    }

    $ws.Range("C2:C26").Select()
}

# ---------------------------------------------------------------------
# quiz6 becomes the active tab/sheet (was metadata_quiz); its own
# selection (A11) is left untouched.
# ---------------------------------------------------------------------
$quiz6 = $wb.Worksheets.Item("quiz6")
$quiz6.Activate()
